# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (column E, rows 16-54) listed the 39 monthly
# periods from 1701 (Jan-2017) through 2003 (Mar-2020) in descending
# order. The refreshed database lists them in ascending (chronological)
# order instead, so every period cell needs to be flipped end-to-end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @(
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}
